# ITSADSSD-21274 - IELTS - Verify Test Components Score
#
# Applies the Config.xlsx changes:
#  - Settings!B2: OrchestratorQueueName value -> IAA_IELTS_Sample_Queue
#  - Settings row 20: new IELTS_URL row (with hyperlink)
#  - Settings!B23:B26: exception email address bpm.ads -> rpa.ads (with hyperlinks)
#  - Settings row 31: allowedItemList row removed
#  - Settings selection moves to A7
#  - Constants row 2: new BusinessProcessName row
#  - Constants selection moves to B37
#  - Assets rows 2-4: new SSO / account rows
#  - Assets selection moves to A3

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")
$assets = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------

# Orchestrator queue name value changed
$settings.Range("B2").Value = "IAA_IELTS_Sample_Queue"

# New row 20 - IELTS website URL
$settings.Range("A20").Value = "IELTS_URL"
$settings.Range("B20").Value = "https://ielts.ucles.org.uk/ielts-trf/"
$settings.Range("C20").Value = "IELTS website login page "

$settings.Range("B20").Style = "Hyperlink"
$settings.Range("B20").IndentLevel = 0

# Exception email address updated from bpm.ads to rpa.ads (DEV/TEST/STAGING/PROD)
$settings.Range("B23").Value = "rpa.ads@its.uq.edu.au"
$settings.Range("B24").Value = "rpa.ads@its.uq.edu.au"
$settings.Range("B25").Value = "rpa.ads@its.uq.edu.au"
$settings.Range("B26").Value = "rpa.ads@its.uq.edu.au"

$settings.Range("B23:B26").Style = "Hyperlink"
$settings.Range("B23:B26").IndentLevel = 0

# allowedItemList row removed entirely
$settings.Range("A31:C31").ClearContents()

# Hyperlinks: B23 (single), B24:B26 (merged mailto range), B20 (IELTS url)
$settings.Hyperlinks.Add($settings.Range("B23"), "mailto:rpa.ads@its.uq.edu.au")
$settings.Hyperlinks.Add($settings.Range("B24:B26"), "mailto:rpa.ads@its.uq.edu.au", [Type]::Missing, [Type]::Missing, "rpa.ads@its.uq.edu.au")
$settings.Hyperlinks.Add($settings.Range("B20"), "https://ielts.ucles.org.uk/ielts-trf/")

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------

$constants.Range("A2").Value = "BusinessProcessName"
$constants.Range("B2").Value = "IAA_IELTS"
$constants.Range("C2").Value = "Business Process Name"

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------

$assets.Range("A2").Value = "SSO_rpaapl01"
$assets.Range("B2").Value = "SSO_rpaapl01"
$assets.Range("C2").Value = "UQ Single Sign On account for Ivy "

$assets.Range("A3").Value = "IAA_IELTS_Account"
$assets.Range("B3").Value = "IAA_IELTS_Account"
$assets.Range("C3").Value = "Ivy login credentails for IELTS website "

$assets.Range("A4").Value = "IAA_IELTS_SPOKED_DB"
$assets.Range("B4").Value = "IAA_IELTS_SPOKED_DB"
$assets.Range("C4").Value = "Credentails to connect with database "

# ---------------------------------------------------------------------------
# Selections - touch Assets and Constants first so Settings ends up as the
# active / tabSelected sheet, matching the original workbook state.
# ---------------------------------------------------------------------------

$assets.Range("A3").Select()
$constants.Range("B37").Select()
$settings.Range("A7").Select()
